$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "Tipo" value in C1 -> it should be "numerico" ---
$ws.Range("C1").Value = "numerico"

# --- Add the new field rows (29-33) describing extra columns added to the dataset ---
# Shared-string table must gain the 5 "name" values first (column A, in row order),
# then the 5 "description" values (column B, in row order), matching the author's edit order.
$newNames = @("region", "latitud_entero", "longitud_entero", "coordenadas", "paises")
$newDescriptions = @(
    "Zona geografica donde se produce el terremoto, formato corto",
    "Latitud expresada indicando solo el entero (sin decimales)",
    "Longitud expresada indicando solo el entero (sin decimales)",
    "Columna donde se indica latitud y longitud en valor entero separados por una coma",
    "País donde se produce el terremoto"
)
$newTypes = @("str", "numerico", "numerico", "str", "str")

$startRow = 29
for ($i = 0; $i -lt $newNames.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $newNames[$i]
}
for ($i = 0; $i -lt $newDescriptions.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $newDescriptions[$i]
}
for ($i = 0; $i -lt $newTypes.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $newTypes[$i]
}

# --- Copy the formatting (yellow fill, per the existing block) onto the new rows ---
$ws.Range("A22:C22").Copy() | Out-Null
$ws.Range("A29:C33").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Update the selected cell to match the author's final cursor position ---
$ws.Range("B35").Select() | Out-Null
